$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell address -> new literal text value, taken from the
# updated "symbol list" snapshot (price + 1h volume % columns).
$updates = [ordered]@{
    "D2" = "288.91"
    "E2" = "-6.84%"
    "D3" = "40.00"
    "E3" = "-2.47%"
    "E4" = "-3.79%"
    "D5" = "0.07361"
    "E5" = "-4.20%"
    "D6" = "4.285"
    "E6" = "-0.13%"
    "D7" = "1.556"
    "E7" = "-9.28%"
    "D8" = "0.9110"
    "E8" = "-2.20%"
    "D9" = "0.1193"
    "E9" = "-6.93%"
    "D10" = "0.1733"
    "E10" = "-5.30%"
    "D11" = "0.08717"
    "E11" = "-4.44%"
    "D12" = "0.04178"
    "E12" = "-0.94%"
    "D13" = "0.1052"
    "E13" = "0.02%"
    "D14" = "0.001278"
    "E14" = "-0.10%"
    "D15" = "0.005961"
    "E15" = "1.21%"
    "D16" = "3.387"
    "E16" = "1.02%"
    "D18" = "0.3301"
    "E18" = "-0.57%"
    "D19" = "7.535"
    "E19" = "-1.94%"
    "E20" = "-0.10%"
    "E21" = "5.95%"
    "E22" = "-4.68%"
    "D23" = "0.001272"
    "E23" = "0.41%"
    "D24" = "0.003887"
    "E24" = "-5.31%"
    "D25" = "0.0001282"
    "E25" = "0.86%"
    "D26" = "0.0003728"
    "D38" = "0.02329"
    "E38" = "-8.56%"
    "D39" = "0.05019"
    "E39" = "-5.47%"
    "D40" = "0.007678"
    "E40" = "-1.98%"
    "E41" = "148.81%"
    "D42" = "0.1272"
    "E42" = "-2.97%"
    "D43" = "0.007383"
    "E43" = "10.96%"
    "D44" = "0.007688"
    "E44" = "-5.14%"
    "D45" = "0.3138"
    "E45" = "1.76%"
    "D46" = "0.00006505"
    "E46" = "-4.15%"
    "E47" = "0.06%"
    "E48" = "11.74%"
    "D49" = "0.004206"
    "E49" = "35.53%"
    "D50" = "0.00002103"
    "E50" = "0.06%"
    "D51" = "0.0002003"
    "E51" = "0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage first so numeric-looking strings (e.g. "40.00",
    # "-6.84%") are kept verbatim instead of being parsed into numbers,
    # which would silently drop trailing zeros / rescale percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset back to the default style so the cell does not end up with a
    # residual text-format style index (the source cells carry no style).
    $cell.Style = "Normal"
}
